$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update destination text values (Column B)
$ws.Range("B2").Value = "Orlando, Florida, United States of America"
$ws.Range("B3").Value = "Phoenix, Arizona, United States of America"

# Update Budget column (Column C) from numeric to descriptive text ranges,
# interleaved with Pool column updates to match original authoring order
$ws.Range("C2").Value = "$300 to $499"

# Update Star Rating column (Column F)
$ws.Range("F2").Value = 4
$ws.Range("F3").Value = 3

# Update Pool column (Column I) from boolean to Yes/No text
$ws.Range("I2").Value = "Yes"
$ws.Range("I3").Value = "No"

$ws.Range("C3").Value = "$100 to $299"

# Update date columns number format (D:E) to mm/dd/yyyy;@
$ws.Range("D2:E3").NumberFormat = "mm/dd/yyyy;@"

# Set column widths for D:E (auto best-fit sized to the new date format, ~10.5546875)
$ws.Range("D1:E1").ColumnWidth = 9.7

# Update active selection cell
$ws.Range("G9").Select()

# Set page orientation to portrait
$ws.PageSetup.Orientation = 1
